# Insert a new weekly data point at row 69 ("Poroto granado" price series),
# pushing the existing rows 69-113 down to 70-114.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(69).Insert()

$ws.Cells.Item(69, 1).Value = 7
$ws.Cells.Item(69, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(69, 3).Value = "Ñuble"
$ws.Cells.Item(69, 4).Value = 44957
$ws.Cells.Item(69, 5).Value = 16
$ws.Cells.Item(69, 6).Value = 100112030
$ws.Cells.Item(69, 7).Value = "Poroto granado"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 30
$ws.Cells.Item(69, 11).Value = 40000
$ws.Cells.Item(69, 12).Value = 40000
$ws.Cells.Item(69, 13).Value = 40000
$ws.Cells.Item(69, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(69, 15).Value = "Región del Maule"
$ws.Cells.Item(69, 16).Value = 1600
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"
